$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 ---
$ws.Range("A9").Value2 = "FewRoles"
$ws.Range("B9").Value2 = 0.65
$ws.Range("C9").Value2 = 1.01
$ws.Range("D9").Value2 = 3.04
$ws.Range("E9").Value2 = 1.01
$ws.Range("F9").Value2 = 0.89
$ws.Range("G9").Value2 = 1.01
$ws.Range("H9").Value2 = 0.89
$ws.Range("I9").Value2 = 1.01
$ws.Range("J9").Value2 = 1.01
$ws.Range("K9").Value2 = 1.01
$ws.Range("L9").Value2 = 1.01
$ws.Range("M9").Value2 = 8

# --- Row 10 ---
$ws.Range("A10").Value2 = "GenericsAndSpecifics"
$ws.Range("B10").Value2 = 0.65
$ws.Range("C10").Value2 = 1.01
$ws.Range("D10").Value2 = 3.04
$ws.Range("E10").Value2 = 1.01
$ws.Range("F10").Value2 = 0.89
$ws.Range("G10").Value2 = 1.01
$ws.Range("H10").Value2 = 0.89
$ws.Range("I10").Value2 = 1.01
$ws.Range("J10").Value2 = 1.01
$ws.Range("K10").Value2 = 1.01
$ws.Range("L10").Value2 = 1.01
$ws.Range("M10").Value2 = 9

# --- Row 11 ---
$ws.Range("A11").Value2 = "SubclassPerTeam"
$ws.Range("B11").Value2 = 1.29
$ws.Range("C11").Value2 = 1.01
$ws.Range("D11").Value2 = 3.04
$ws.Range("E11").Value2 = 1.01
$ws.Range("F11").Value2 = 0.89
$ws.Range("G11").Value2 = 1.01
$ws.Range("H11").Value2 = 2.67
$ws.Range("I11").Value2 = 1.01
$ws.Range("J11").Value2 = 1.01
$ws.Range("K11").Value2 = 1.01
$ws.Range("L11").Value2 = 1.01
$ws.Range("M11").Value2 = 10

# --- Row 12 ---
$ws.Range("A12").Value2 = "HierarchyOfFactories"
$ws.Range("B12").Value2 = 0.65
$ws.Range("C12").Value2 = 1.01
$ws.Range("D12").Value2 = 3.04
$ws.Range("E12").Value2 = 1.01
$ws.Range("F12").Value2 = 0.89
$ws.Range("G12").Value2 = 1.01
$ws.Range("H12").Value2 = 0.89
$ws.Range("I12").Value2 = 1.01
$ws.Range("J12").Value2 = 1.01
$ws.Range("K12").Value2 = 1.01
$ws.Range("L12").Value2 = 1.01
$ws.Range("M12").Value2 = 11

# --- Row 13 ---
$ws.Range("A13").Value2 = "LooseInterfaces"
$ws.Range("B13").Value2 = 0.65
$ws.Range("C13").Value2 = 1.01
$ws.Range("D13").Value2 = 3.04
$ws.Range("E13").Value2 = 1.01
$ws.Range("F13").Value2 = 0.89
$ws.Range("G13").Value2 = 1.01
$ws.Range("H13").Value2 = 0.89
$ws.Range("I13").Value2 = 1.01
$ws.Range("J13").Value2 = 1.01
$ws.Range("K13").Value2 = 1.01
$ws.Range("L13").Value2 = 1.01
$ws.Range("M13").Value2 = 12

# --- Row 14 ---
$ws.Range("A14").Value2 = "FewRoles"
$ws.Range("B14").Value2 = 0.65
$ws.Range("C14").Value2 = 1.01
$ws.Range("D14").Value2 = 3.04
$ws.Range("E14").Value2 = 1.01
$ws.Range("F14").Value2 = 0.89
$ws.Range("G14").Value2 = 1.01
$ws.Range("H14").Value2 = 0.89
$ws.Range("I14").Value2 = 1.01
$ws.Range("J14").Value2 = 1.01
$ws.Range("K14").Value2 = 1.01
$ws.Range("L14").Value2 = 1.01
$ws.Range("M14").Value2 = 13

# --- Row 15 ---
$ws.Range("A15").Value2 = "ParserBuilder"
$ws.Range("B15").Value2 = 0.65
$ws.Range("C15").Value2 = 1.01
$ws.Range("D15").Value2 = 3.04
$ws.Range("E15").Value2 = 1.01
$ws.Range("F15").Value2 = 0.89
$ws.Range("G15").Value2 = 1.01
$ws.Range("H15").Value2 = 0.89
$ws.Range("I15").Value2 = 1.01
$ws.Range("J15").Value2 = 1.01
$ws.Range("K15").Value2 = 1.01
$ws.Range("L15").Value2 = 1.01
$ws.Range("M15").Value2 = 14

# --- Row 16 ---
$ws.Range("A16").Value2 = "ArchitectControlsProduct"
$ws.Range("B16").Value2 = 0.65
$ws.Range("C16").Value2 = 1.01
$ws.Range("D16").Value2 = 3.04
$ws.Range("E16").Value2 = 1.01
$ws.Range("F16").Value2 = 0.89
$ws.Range("G16").Value2 = 1.01
$ws.Range("H16").Value2 = 0.89
$ws.Range("I16").Value2 = 1.01
$ws.Range("J16").Value2 = 1.01
$ws.Range("K16").Value2 = 1.01
$ws.Range("L16").Value2 = 1.01
$ws.Range("M16").Value2 = 15

# --- Row 17 ---
$ws.Range("A17").Value2 = "DeployAlongTheGrain"
$ws.Range("B17").Value2 = 0.65
$ws.Range("C17").Value2 = 1.01
$ws.Range("D17").Value2 = 3.04
$ws.Range("E17").Value2 = 1.01
$ws.Range("F17").Value2 = 0.89
$ws.Range("G17").Value2 = 1.01
$ws.Range("H17").Value2 = 0.89
$ws.Range("I17").Value2 = 1.01
$ws.Range("J17").Value2 = 1.01
$ws.Range("K17").Value2 = 1.01
$ws.Range("L17").Value2 = 1.01
$ws.Range("M17").Value2 = 16

# --- Row 18 ---
$ws.Range("A18").Value2 = "DeveloperControlsProcess"
$ws.Range("B18").Value2 = 0.65
$ws.Range("C18").Value2 = 1.01
$ws.Range("D18").Value2 = 3.04
$ws.Range("E18").Value2 = 1.01
$ws.Range("F18").Value2 = 0.89
$ws.Range("G18").Value2 = 1.01
$ws.Range("H18").Value2 = 0.89
$ws.Range("I18").Value2 = 1.01
$ws.Range("J18").Value2 = 1.01
$ws.Range("K18").Value2 = 1.01
$ws.Range("L18").Value2 = 1.01
$ws.Range("M18").Value2 = 17

# --- Row 19 ---
$ws.Range("A19").Value2 = "ArchitectureTeam"
$ws.Range("B19").Value2 = 0.65
$ws.Range("C19").Value2 = 1.01
$ws.Range("D19").Value2 = 3.04
$ws.Range("E19").Value2 = 1.01
$ws.Range("F19").Value2 = 0.89
$ws.Range("G19").Value2 = 1.01
$ws.Range("H19").Value2 = 0.89
$ws.Range("I19").Value2 = 1.01
$ws.Range("J19").Value2 = 1.01
$ws.Range("K19").Value2 = 1.01
$ws.Range("L19").Value2 = 1.01
$ws.Range("M19").Value2 = 18

# Re-apply the bordered "highlight" style (same one used on A2/M2/M4/M6/M8)
# to the specific cells that keep it in the new rows.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("M2").Copy()
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("D18").Select()
